$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 4018
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 6022.5
$ws.Range("K16").Value = 9
$ws.Range("L16").Value = 6022.5
$ws.Range("M16").Value = 221
$ws.Range("N16").Value = -6482.5
$ws.Range("H17").Value = 1980.6666
$ws.Range("J17").Value = 1980.6666
$ws.Range("L17").Value = 5941.9998
$ws.Range("N17").Value = -6277.9998
$ws.Range("H31").Value = 105000
$ws.Range("I31").Value = 105000
$ws.Range("K31").Value = 315000
$ws.Range("M31").Value = -314770
$ws.Range("H113").Value = 3232.6458
$ws.Range("J113").Value = 3082.4482
$ws.Range("L113").Value = 3082.4482
$ws.Range("N113").Value = -9590.448199999999
$ws.Range("H116").Value = 2928.5715
$ws.Range("I116").Value = 3274.5
$ws.Range("K116").Value = 3274.5
$ws.Range("M116").Value = 167.5
$ws.Range("H137").Value = 2587.1904
$ws.Range("I137").Value = 2593.5
$ws.Range("J137").Value = 2549.3333
$ws.Range("K137").Value = 7780.5
$ws.Range("L137").Value = 7647.999899999999
$ws.Range("M137").Value = -5230.5
$ws.Range("N137").Value = -12747.9999
$ws.Range("H138").Value = 2411.561
$ws.Range("I138").Value = 1069.8235
$ws.Range("J138").Value = 3361.9583
$ws.Range("K138").Value = 3209.4705
$ws.Range("L138").Value = 10085.8749
$ws.Range("M138").Value = 1930.5295
$ws.Range("N138").Value = -20365.8749
$ws.Range("H141").Value = 4136.467
$ws.Range("I141").Value = 4136.467
$ws.Range("K141").Value = 12409.401
$ws.Range("M141").Value = -7229.400999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 17260.45
$ws.Range("I2").Value = 27506.25
$ws.Range("K2").Value = 27506.25
$ws.Range("M2").Value = -27393.25
$ws.Range("H32").Value = 9145.833000000001
$ws.Range("I32").Value = 8261.75
$ws.Range("K32").Value = 8261.75
$ws.Range("M32").Value = -7974.75
$ws.Range("H102").Value = 629.7143
$ws.Range("I102").Value = 629.7143
$ws.Range("K102").Value = 629.7143
$ws.Range("M102").Value = 992.2857
$ws.Range("H116").Value = 17260.45
$ws.Range("I116").Value = 27506.25
$ws.Range("K116").Value = 27506.25
$ws.Range("M116").Value = -25212.25
$ws.Range("H132").Value = 5696.353
$ws.Range("I132").Value = 5574.4
$ws.Range("J132").Value = 5870.5713
$ws.Range("K132").Value = 16723.2
$ws.Range("L132").Value = 17611.7139
$ws.Range("M132").Value = -14193.2
$ws.Range("N132").Value = -22671.7139

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 17260.45
$ws.Range("I3").Value = 27506.25
$ws.Range("K3").Value = 27506.25
$ws.Range("M3").Value = -27392.25
$ws.Range("H86").Value = 2910.7666
$ws.Range("I86").Value = 1515
$ws.Range("J86").Value = 4505.9287
$ws.Range("K86").Value = 1515
$ws.Range("L86").Value = 4505.9287
$ws.Range("M86").Value = -392
$ws.Range("N86").Value = -6751.9287
$ws.Range("H89").Value = 2910.7666
$ws.Range("I89").Value = 1515
$ws.Range("J89").Value = 4505.9287
$ws.Range("K89").Value = 7575
$ws.Range("L89").Value = 22529.6435
$ws.Range("M89").Value = -1959
$ws.Range("N89").Value = -33761.64350000001
$ws.Range("H99").Value = 20028.545
$ws.Range("I99").Value = 23574.834
$ws.Range("K99").Value = 23574.834
$ws.Range("M99").Value = -22076.834
$ws.Range("H105").Value = 1876.4
$ws.Range("I105").Value = 1876.4
$ws.Range("K105").Value = 1876.4
$ws.Range("M105").Value = -129.4000000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4320132
$ws.Range("J99").Value = 5719449
$ws.Range("L99").Value = 5719449
$ws.Range("N99").Value = -5722445
$ws.Range("H126").Value = 4320132
$ws.Range("J126").Value = 5719449
$ws.Range("L126").Value = 17158347
$ws.Range("N126").Value = -17163287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 4997.25
$ws.Range("J103").Value = 5000
$ws.Range("L103").Value = 15000
$ws.Range("N103").Value = -16758
$ws.Range("H129").Value = 2428.923
$ws.Range("J129").Value = 2729.6667
$ws.Range("L129").Value = 8189.000100000001
$ws.Range("N129").Value = -18189.0001
$ws.Range("H138").Value = 3280
$ws.Range("I138").Value = 2690.25
$ws.Range("J138").Value = 7998
$ws.Range("K138").Value = 8070.75
$ws.Range("L138").Value = 23994
$ws.Range("M138").Value = -2930.75
$ws.Range("N138").Value = -34274
$ws.Range("H139").Value = 30308530
$ws.Range("I139").Value = 52634500
$ws.Range("K139").Value = 157903500
$ws.Range("M139").Value = -157898360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 72321.625
$ws.Range("I80").Value = 103508.82
$ws.Range("J80").Value = 3709.8
$ws.Range("K80").Value = 103508.82
$ws.Range("L80").Value = 3709.8
$ws.Range("M80").Value = -102510.82
$ws.Range("N80").Value = -5705.8
$ws.Range("H83").Value = 72321.625
$ws.Range("I83").Value = 103508.82
$ws.Range("J83").Value = 3709.8
$ws.Range("K83").Value = 517544.1
$ws.Range("L83").Value = 18549
$ws.Range("M83").Value = -512552.1
$ws.Range("N83").Value = -28533
$ws.Range("H107").Value = 305.9
$ws.Range("J107").Value = 174.75
$ws.Range("L107").Value = 174.75
$ws.Range("N107").Value = -4014.75
$ws.Range("H122").Value = 4249.636
$ws.Range("I122").Value = 4746.1333
$ws.Range("K122").Value = 14238.3999
$ws.Range("M122").Value = -11788.3999
$ws.Range("H123").Value = 37500
$ws.Range("J123").Value = 37500
$ws.Range("L123").Value = 37500
$ws.Range("N123").Value = -42400
$ws.Range("H132").Value = 8257.299999999999
$ws.Range("I132").Value = 8223.25
$ws.Range("K132").Value = 24669.75
$ws.Range("M132").Value = -22139.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3596.3
$ws.Range("I7").Value = 2471.6
$ws.Range("J7").Value = 4721
$ws.Range("K7").Value = 2471.6
$ws.Range("L7").Value = 4721
$ws.Range("M7").Value = -2359.6
$ws.Range("N7").Value = -4945
$ws.Range("H13").Value = 8876
$ws.Range("I13").Value = 4696.25
$ws.Range("J13").Value = 11662.5
$ws.Range("K13").Value = 4696.25
$ws.Range("L13").Value = 11662.5
$ws.Range("M13").Value = -4556.25
$ws.Range("N13").Value = -11942.5
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -705
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -893
$ws.Range("H40").Value = 8456.069
$ws.Range("I40").Value = 9214.467000000001
$ws.Range("K40").Value = 9214.467000000001
$ws.Range("M40").Value = -9078.467000000001
$ws.Range("H122").Value = 5110.731
$ws.Range("I122").Value = 3790.7368
$ws.Range("K122").Value = 11372.2104
$ws.Range("M122").Value = -8922.2104
$ws.Range("H126").Value = 3596.3
$ws.Range("I126").Value = 2471.6
$ws.Range("J126").Value = 4721
$ws.Range("K126").Value = 7414.799999999999
$ws.Range("L126").Value = 14163
$ws.Range("M126").Value = -4944.799999999999
$ws.Range("N126").Value = -19103
$ws.Range("H136").Value = 4596.7354
$ws.Range("I136").Value = 2695.4666
$ws.Range("K136").Value = 8086.399800000001
$ws.Range("M136").Value = -5536.399800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 64265.332
$ws.Range("J75").Value = 64265.332
$ws.Range("L75").Value = 64265.332
$ws.Range("N75").Value = -66137.33199999999
$ws.Range("H78").Value = 64265.332
$ws.Range("J78").Value = 64265.332
$ws.Range("L78").Value = 192795.996
$ws.Range("N78").Value = -202155.996
$ws.Range("H81").Value = 5100.5864
$ws.Range("I81").Value = 5886.7827
$ws.Range("K81").Value = 11773.5654
$ws.Range("M81").Value = -10712.5654
$ws.Range("H84").Value = 5100.5864
$ws.Range("I84").Value = 5886.7827
$ws.Range("K84").Value = 58867.827
$ws.Range("M84").Value = -53563.827
$ws.Range("H122").Value = 4334.5713
$ws.Range("I122").Value = 1855.7407
$ws.Range("J122").Value = 12700.625
$ws.Range("K122").Value = 5567.2221
$ws.Range("L122").Value = 38101.875
$ws.Range("M122").Value = -3117.2221
$ws.Range("N122").Value = -43001.875
$ws.Range("H126").Value = 2673.3333
$ws.Range("I126").Value = 2258
$ws.Range("J126").Value = 4750
$ws.Range("K126").Value = 6774
$ws.Range("L126").Value = 14250
$ws.Range("M126").Value = -4304
$ws.Range("N126").Value = -19190
$ws.Range("H132").Value = 3165.5454
$ws.Range("I132").Value = 2699.8462
$ws.Range("K132").Value = 8099.5386
$ws.Range("M132").Value = -5569.5386

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N22").ClearContents()
$ws.Range("N27").ClearContents()
